$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.003.38"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.15"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.10"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5092"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06363"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.296"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.70"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5485"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7747"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.41"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.024.03"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.02"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.464"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.970"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.134"
$ws.Range("E22").Value = "  +1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.93"
$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1262"
$ws.Range("E26").Value = "  +10.32%  "

$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.67"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.242"
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04889"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.283"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.220"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.553"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.377"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9198"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.573"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5552"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.104.27"
$ws.Range("E38").Value = "  -2.29%  "

$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.622"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8057"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.96"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("E44").Value = "  -4.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.778.96"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.40"
$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05188"
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.570"
$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("E51").Value = "  -0.32%  "
